$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Set Cycle" column (D3:D9) to reflect new currency-db based
# settlement cycle values: "T+2" for rows 3-8, "T+3" for row 9.
$ws.Range("D3").Value = "T+2"
$ws.Range("D4").Value = "T+2"
$ws.Range("D5").Value = "T+2"
$ws.Range("D6").Value = "T+2"
$ws.Range("D7").Value = "T+2"
$ws.Range("D8").Value = "T+2"
$ws.Range("D9").Value = "T+3"

# Move the active selection to D10 (matches the saved view state in the
# edited workbook).
$ws.Range("D10").Select()
